# Update Betfair Back/Lay odds for 2026-01-15 workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 3.9
$ws.Range("T2").Value = 1.96

# Row 3
$ws.Range("F3").Value = 1.22
$ws.Range("K3").Value = 9
$ws.Range("N3").Value = 8
$ws.Range("P3").Value = 3.8
$ws.Range("Q3").Value = 1.3
$ws.Range("R3").Value = 2.12
$ws.Range("S3").Value = 1.75
$ws.Range("U3").Value = 2.1
$ws.Range("X3").Value = 60
$ws.Range("Y3").Value = 70
$ws.Range("Z3").Value = 170
$ws.Range("AA3").Value = 490
$ws.Range("AB3").Value = 17.5
$ws.Range("AC3").Value = 24
$ws.Range("AE3").Value = 190
$ws.Range("AF3").Value = 14
$ws.Range("AG3").Value = 15.5
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 14.5
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 2.94
$ws.Range("AO3").Value = 150

# Row 4
$ws.Range("H4").Value = 2.24
$ws.Range("I4").Value = 2.26
$ws.Range("O4").Value = 1.47
$ws.Range("Q4").Value = 2.4
$ws.Range("T4").Value = 2.04

# Row 5
$ws.Range("J5").Value = 3.3
$ws.Range("U5").Value = 1.97
$ws.Range("Z5").Value = 18.5

# Row 6
$ws.Range("F6").Value = 3.2
$ws.Range("G6").Value = 3.25
$ws.Range("P6").Value = 1.83
$ws.Range("S6").Value = 3.95
$ws.Range("AA6").Value = 38
